# Generate Report for Handoff
# Appends two new "Ready for handoff" rows (for f5b8b817-... and f785546e-...)
# to the Overview / zh-cn / de-de sheets, wires up hyperlinks for the new
# source-file cells, and grows each sheet's table to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1) Overview sheet - two new rows (6 and 7)
# ---------------------------------------------------------------------

$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = "'"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "'2016-09-01 18:48:38"

$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = "'"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "'2016-09-01 18:48:38"

$wsOverview.Range("A6").Value = "f5b8b817-171a-4d9a-9c10-318d86d2d96a.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5b8b817171a4d9a9c10318d86d2d96aaaaaaaa/e2e/f5b8b817-171a-4d9a-9c10-318d86d2d96a.md", "", "", "e2e\f5b8b817-171a-4d9a-9c10-318d86d2d96a.md")

$wsOverview.Range("A7").Value = "f785546e-6d87-42af-8ccc-5203ee0dac32.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f785546e6d8742af8ccc5203ee0dac32aaaaaaaa/e2e/f785546e-6d87-42af-8ccc-5203ee0dac32.md", "", "", "e2e\f785546e-6d87-42af-8ccc-5203ee0dac32.md")

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------
# 2) zh-cn sheet - two new rows (6 and 7)
# ---------------------------------------------------------------------

$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "'False"
$wsZhCn.Range("G6").Value = "f5b8b817-171a-4d9a-9c10-318d86d2d96a.a1e932613a9dc9f847f31387e8a65a56127dcec8.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "'2016-09-01 18:48:33"
$wsZhCn.Range("I6").Value = "'"
$wsZhCn.Range("J6").Value = "'"
$wsZhCn.Range("K6").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L6").Value = "'"
$wsZhCn.Range("M6").Value = "'True"
$wsZhCn.Range("N6").Value = "'"
$wsZhCn.Range("O6").Value = "'False"
$wsZhCn.Range("P6").Value = "'"

$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "'False"
$wsZhCn.Range("G7").Value = "f785546e-6d87-42af-8ccc-5203ee0dac32.cdfaf3490c42c8626702b019a8d20abff4381555.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "'2016-09-01 18:48:33"
$wsZhCn.Range("I7").Value = "'"
$wsZhCn.Range("J7").Value = "'"
$wsZhCn.Range("K7").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L7").Value = "'"
$wsZhCn.Range("M7").Value = "'True"
$wsZhCn.Range("N7").Value = "'"
$wsZhCn.Range("O7").Value = "'False"
$wsZhCn.Range("P7").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5b8b817171a4d9a9c10318d86d2d96aaaaaaaa/e2e/f5b8b817-171a-4d9a-9c10-318d86d2d96a.md", "", "", "f5b8b817-171a-4d9a-9c10-318d86d2d96a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f785546e6d8742af8ccc5203ee0dac32aaaaaaaa/e2e/f785546e-6d87-42af-8ccc-5203ee0dac32.md", "", "", "f785546e-6d87-42af-8ccc-5203ee0dac32.md")

$loZhCn = $wsZhCn.ListObjects.Item("zh_cn")
$loZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------
# 3) de-de sheet - two new rows (6 and 7)
# ---------------------------------------------------------------------

$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "'False"
$wsDeDe.Range("G6").Value = "f5b8b817-171a-4d9a-9c10-318d86d2d96a.a1e932613a9dc9f847f31387e8a65a56127dcec8.de-de.xlf"
$wsDeDe.Range("H6").Value = "'2016-09-01 18:48:38"
$wsDeDe.Range("I6").Value = "'"
$wsDeDe.Range("J6").Value = "'"
$wsDeDe.Range("K6").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L6").Value = "'"
$wsDeDe.Range("M6").Value = "'True"
$wsDeDe.Range("N6").Value = "'"
$wsDeDe.Range("O6").Value = "'False"
$wsDeDe.Range("P6").Value = "'"

$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "'False"
$wsDeDe.Range("G7").Value = "f785546e-6d87-42af-8ccc-5203ee0dac32.cdfaf3490c42c8626702b019a8d20abff4381555.de-de.xlf"
$wsDeDe.Range("H7").Value = "'2016-09-01 18:48:38"
$wsDeDe.Range("I7").Value = "'"
$wsDeDe.Range("J7").Value = "'"
$wsDeDe.Range("K7").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L7").Value = "'"
$wsDeDe.Range("M7").Value = "'True"
$wsDeDe.Range("N7").Value = "'"
$wsDeDe.Range("O7").Value = "'False"
$wsDeDe.Range("P7").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5b8b817171a4d9a9c10318d86d2d96aaaaaaaa/e2e/f5b8b817-171a-4d9a-9c10-318d86d2d96a.md", "", "", "f5b8b817-171a-4d9a-9c10-318d86d2d96a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f785546e6d8742af8ccc5203ee0dac32aaaaaaaa/e2e/f785546e-6d87-42af-8ccc-5203ee0dac32.md", "", "", "f785546e-6d87-42af-8ccc-5203ee0dac32.md")

$loDeDe = $wsDeDe.ListObjects.Item("de_de")
$loDeDe.Resize($wsDeDe.Range("A1:P7"))

Write-Host "Report regenerated: Overview/zh-cn/de-de now include f5b8b817... and f785546e... handoff rows."
